$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.3917
$ws.Range("C6").Value = -12.141
$ws.Range("C7").Value = -12.95209999999999
$ws.Range("C8").Value = -12.4149
$ws.Range("C16").Value = -14.41379999999999
$ws.Range("C20").Value = -11.6442
$ws.Range("C21").Value = -12.3185
$ws.Range("C28").Value = -12.5847
$ws.Range("C29").Value = -11.47870000000001
$ws.Range("C30").Value = -12.11409999999999
$ws.Range("C32").Value = -12.4386
$ws.Range("C40").Value = -12.9552
$ws.Range("C46").Value = -14.66309999999999
$ws.Range("C51").Value = -11.3471
$ws.Range("C52").Value = -11.2666
$ws.Range("C57").Value = -14.34399999999999
$ws.Range("C59").Value = -12.82579999999999
$ws.Range("C62").Value = -14.38339999999999
$ws.Range("C66").Value = -11.0761
$ws.Range("C73").Value = -12.7569
$ws.Range("C74").Value = -11.87940000000001
$ws.Range("C77").Value = -12.0825
$ws.Range("C92").Value = -11.2422
$ws.Range("C100").Value = -12.93439999999999
